$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate rows 26-101 (columns B:F) with recorded tick data.
# Data layout per row: [row, Date(B), Time(C), Price(D), Volume(E)]; F is always "-".
$rows = @(
    @(26, 44488, 0.78471064814814817, 4.9809999999999999, 1),
    @(27, 44488, 0.78465277777777775, 4.9850000000000003, 1),
    @(28, 44488, 0.78396990740740735, 4.9850000000000003, 1),
    @(29, 44488, 0.78379629629629621, 4.9850000000000003, 1),
    @(30, 44488, 0.78361111111111104, 4.9829999999999997, 5),
    @(31, 44488, 0.78297453703703701, 4.9809999999999999, 1),
    @(32, 44488, 0.7819328703703704, 4.9800000000000004, 1),
    @(33, 44488, 0.781712962962963, 4.9809999999999999, 1),
    @(34, 44488, 0.78123842592592585, 4.9809999999999999, 1),
    @(35, 44488, 0.78118055555555566, 4.984, 1),
    @(36, 44488, 0.780787037037037, 4.9809999999999999, 1),
    @(37, 44488, 0.78008101851851863, 4.9820000000000002, 1),
    @(38, 44488, 0.77983796296296293, 4.984, 1),
    @(39, 44488, 0.77938657407407408, 4.984, 2),
    @(40, 44488, 0.77915509259259252, 4.9809999999999999, 1),
    @(41, 44488, 0.77913194444444445, 4.984, 1),
    @(42, 44488, 0.7785185185185185, 4.9829999999999997, 1),
    @(43, 44488, 0.77798611111111116, 4.9829999999999997, 1),
    @(44, 44488, 0.77792824074074074, 4.9829999999999997, 1),
    @(45, 44488, 0.77789351851851851, 4.9829999999999997, 1),
    @(46, 44488, 0.7777546296296296, 4.9809999999999999, 2),
    @(47, 44488, 0.7775347222222222, 4.9800000000000004, 1),
    @(48, 44488, 0.77696759259259263, 4.9770000000000003, 1),
    @(49, 44488, 0.77695601851851848, 4.9800000000000004, 1),
    @(50, 44488, 0.77637731481481476, 4.976, 1),
    @(51, 44488, 0.77568287037037031, 4.9800000000000004, 1),
    @(52, 44488, 0.77568287037037031, 4.9740000000000002, 1),
    @(53, 44488, 0.77545138888888887, 4.9740000000000002, 1),
    @(54, 44488, 0.77542824074074079, 4.976, 1),
    @(55, 44488, 0.77541666666666664, 4.9770000000000003, 1),
    @(56, 44488, 0.77541666666666664, 4.9779999999999998, 3),
    @(57, 44488, 0.77541666666666664, 4.9779999999999998, 1),
    @(58, 44488, 0.77533564814814815, 4.9820000000000002, 1),
    @(59, 44488, 0.77483796296296292, 4.9790000000000001, 1),
    @(60, 44488, 0.77482638888888899, 4.9790000000000001, 1),
    @(61, 44488, 0.7747222222222222, 4.9790000000000001, 3),
    @(62, 44488, 0.77432870370370377, 4.9790000000000001, 1),
    @(63, 44488, 0.77429398148148154, 4.9790000000000001, 1),
    @(64, 44488, 0.77302083333333327, 4.9800000000000004, 1),
    @(65, 44488, 0.77293981481481477, 4.9809999999999999, 1),
    @(66, 44488, 0.77293981481481477, 4.9820000000000002, 1),
    @(67, 44488, 0.7726736111111111, 4.9829999999999997, 1),
    @(68, 44488, 0.77229166666666671, 4.9829999999999997, 1),
    @(69, 44488, 0.77217592592592599, 4.9820000000000002, 2),
    @(70, 44488, 0.77186342592592594, 4.9809999999999999, 1),
    @(71, 44488, 0.77164351851851853, 4.9800000000000004, 3),
    @(72, 44488, 0.77164351851851853, 4.9800000000000004, 1),
    @(73, 44488, 0.77144675925925921, 4.9800000000000004, 1),
    @(74, 44488, 0.77127314814814818, 4.9770000000000003, 2),
    @(75, 44488, 0.77087962962962964, 4.9770000000000003, 2),
    @(76, 44488, 0.77086805555555549, 4.9790000000000001, 3),
    @(77, 44488, 0.77008101851851851, 4.9770000000000003, 1),
    @(78, 44488, 0.77005787037037043, 4.9770000000000003, 3),
    @(79, 44488, 0.76995370370370375, 4.976, 2),
    @(80, 44488, 0.76995370370370375, 4.976, 1),
    @(81, 44488, 0.76861111111111102, 4.976, 2),
    @(82, 44488, 0.76848379629629626, 4.9749999999999996, 1),
    @(83, 44488, 0.76806712962962964, 4.9729999999999999, 1),
    @(84, 44488, 0.76804398148148145, 4.9749999999999996, 1),
    @(85, 44488, 0.76774305555555555, 4.9729999999999999, 1),
    @(86, 44488, 0.76760416666666664, 4.9729999999999999, 1),
    @(87, 44488, 0.76759259259259249, 4.9729999999999999, 1),
    @(88, 44488, 0.76759259259259249, 4.9729999999999999, 1),
    @(89, 44488, 0.76759259259259249, 4.9729999999999999, 1),
    @(90, 44488, 0.76759259259259249, 4.9770000000000003, 1),
    @(91, 44488, 0.76759259259259249, 4.9749999999999996, 2),
    @(92, 44488, 0.76758101851851857, 4.976, 3),
    @(93, 44488, 0.76731481481481489, 4.976, 1),
    @(94, 44488, 0.76729166666666659, 4.976, 4),
    @(95, 44488, 0.76714120370370376, 4.976, 1),
    @(96, 44488, 0.76695601851851858, 4.9770000000000003, 1),
    @(97, 44488, 0.76688657407407401, 4.9770000000000003, 1),
    @(98, 44488, 0.76677083333333329, 4.9770000000000003, 2),
    @(99, 44488, 0.76675925925925925, 4.9770000000000003, 1),
    @(100, 44488, 0.76674768518518521, 4.9790000000000001, 1),
    @(101, 44488, 0.76666666666666661, 4.9790000000000001, 1)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]   # B: Date
    $ws.Cells.Item($r, 3).Value = $row[2]   # C: Time
    $ws.Cells.Item($r, 4).Value = $row[3]   # D: Price
    $ws.Cells.Item($r, 5).Value = $row[4]   # E: Volume
    $ws.Cells.Item($r, 6).Value = "-"       # F: Symbol
}

# Reflect the just-filled range as the active selection.
$ws.Range("B1:F101").Select() | Out-Null
